# Y2_B2526_Blood_&_lymphatics_schedule.xlsx
# Update "Parasitology SGD/POS" session duration (column G) from 75 to 90
# minutes across the schedule, and (re)apply the AutoFilter over the used
# range A1:G154 (uploaded via the attendance app).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Duration (column G) changes from 75 to 90.
$rows = @(17, 18, 34, 35, 51, 52, 68, 69, 85, 86, 102, 103, 119, 120, 136, 137)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = 90
}

# Re-apply the worksheet AutoFilter over the full data range, which
# registers the hidden _xlnm._FilterDatabase defined name scoped to this
# sheet (as seen in the saved workbook.xml).
$nm = $ws.Names.Add("_xlnm._FilterDatabase", "='Blood_&_lymphatics'!`$A`$1:`$G`$154")
$nm.Visible = $false
